$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $savedStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "43.803.30"
$ws.Range("E2").Value = "  +0.00%  "
Set-TextValue $ws.Range("D3") "2.255.82"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue $ws.Range("D5") "230.42"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +2.66%  "
Set-TextValue $ws.Range("D7") "64.16"
$ws.Range("E7").Value = "  +4.92%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.449"
$ws.Range("E9").Value = "  +6.90%  "
Set-TextValue $ws.Range("D10") "0.0974"
$ws.Range("E10").Value = "  +4.43%  "
Set-TextValue $ws.Range("D11") "57.00"
$ws.Range("E11").Value = "  -1.62%  "
Set-TextValue $ws.Range("D12") "26.74"
$ws.Range("E12").Value = "  +13.08%  "
$ws.Range("E13").Value = "  +2.22%  "
Set-TextValue $ws.Range("D14") "2.596.17"
$ws.Range("E14").Value = "  -0.44%  "
Set-TextValue $ws.Range("D15") "15.60"
$ws.Range("E15").Value = "  +0.29%  "
Set-TextValue $ws.Range("D16") "6.09"
$ws.Range("E16").Value = "  +5.44%  "
Set-TextValue $ws.Range("D17") "0.834"
$ws.Range("E17").Value = "  +3.20%  "
Set-TextValue $ws.Range("D18") "2.260.10"
$ws.Range("E18").Value = "  -0.45%  "
Set-TextValue $ws.Range("D19") "43.790.21"
$ws.Range("E19").Value = "  +0.17%  "
Set-TextValue $ws.Range("D20") "0.0₃0987"
$ws.Range("E20").Value = "  +5.54%  "
Set-TextValue $ws.Range("D21") "73.27"
$ws.Range("E21").Value = "  +0.35%  "
Set-TextValue $ws.Range("D22") "6.05"
$ws.Range("E22").Value = "  -2.47%  "
Set-TextValue $ws.Range("D23") "250.44"
$ws.Range("E23").Value = "  -1.05%  "
Set-TextValue $ws.Range("D24") "0.999"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -4.32%  "
Set-TextValue $ws.Range("D26") "2.25"
$ws.Range("E26").Value = "  -8.85%  "
Set-TextValue $ws.Range("D27") "10.04"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +22.91%  "
Set-TextValue $ws.Range("D29") "170.88"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -0.64%  "
Set-TextValue $ws.Range("D31") "20.89"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("E33").Value = "  +2.92%  "
Set-TextValue $ws.Range("D34") "0.0704"
$ws.Range("E34").Value = "  +7.05%  "
Set-TextValue $ws.Range("D35") "4.79"
$ws.Range("E35").Value = "  +0.60%  "
Set-TextValue $ws.Range("D36") "4.89"
$ws.Range("E36").Value = "  -3.49%  "
Set-TextValue $ws.Range("D37") "3.73"
$ws.Range("E37").Value = "  +4.01%  "
Set-TextValue $ws.Range("D38") "6.48"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -4.05%  "
Set-TextValue $ws.Range("D40") "0.0260"
$ws.Range("E40").Value = "  +4.46%  "
$ws.Range("E41").Value = "  +0.04%  "
Set-TextValue $ws.Range("D42") "0.000223"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  +4.04%  "
Set-TextValue $ws.Range("D44") "0.0966"
$ws.Range("E44").Value = "  -1.94%  "
Set-TextValue $ws.Range("D45") "8.20"
$ws.Range("E45").Value = "  -6.13%  "
Set-TextValue $ws.Range("D46") "97.74"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  -0.94%  "
Set-TextValue $ws.Range("D48") "4.41"
$ws.Range("E48").Value = "  -2.13%  "
Set-TextValue $ws.Range("D49") "2.37"
$ws.Range("E49").Value = "  +5.87%  "
Set-TextValue $ws.Range("D50") "10.16"
$ws.Range("E50").Value = "  +6.34%  "
Set-TextValue $ws.Range("D51") "1.436.16"
$ws.Range("E51").Value = "  -2.87%  "

Write-Output "done"